$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'list(au_id = c("https://openalex.org/A5035612138", "https://openalex.org/A5036981754", "https://openalex.org/A5080741438", "https://openalex.org/A5005923405"), au_display_name = c("Peter Lykke Eriksen", "Lars Djernes", "Hendrik Vilstrup", "Peter Ott"), au_orcid = c("https://orcid.org/0000-0003-1115-8308", NA, NA, "https://orcid.org/0000-0002-3088-1983"), author_position = c("first", "middle", "middle", "last"), au_affiliation_raw = c("Department of Hepatology and Gastroenterology, Aarhus University Hospital, Denmark", 
"Department of Anaesthesiology and Intensive Care, Aarhus University Hospital, Denmark; Department of Anaesthesiology and Intensive Care, Viborg Regional Hospital, Denmark", "Department of Clinical Medicine, Aarhus University, Denmark; Department of Hepatology and Gastroenterology, Aarhus University Hospital, Denmark", "Department of Clinical Medicine, Aarhus University, Denmark; Department of Hepatology and Gastroenterology, Aarhus University Hospital, Denmark"), institution_id = c("https://openalex.org/I2802335433", 
"https://openalex.org/I2802335433", "https://openalex.org/I2802335433", "https://openalex.org/I2802335433"), institution_display_name = c("Aarhus University Hospital", "Aarhus University Hospital", "Aarhus University Hospital", "Aarhus University Hospital"), institution_ror = c("https://ror.org/040r8fr65", "https://ror.org/040r8fr65", "https://ror.org/040r8fr65", "https://ror.org/040r8fr65"), institution_country_code = c("DK", "DK", "DK", "DK"), institution_type = c("healthcare", "healthcare", "healthcare", 
"healthcare"), institution_lineage = c("https://openalex.org/I2802335433", "https://openalex.org/I2802335433", "https://openalex.org/I2802335433", "https://openalex.org/I2802335433"))'

$ws.Range("AF2").Value = 'c("https://openalex.org/W3031953493", "https://openalex.org/W2791173279", "https://openalex.org/W2809943007", "https://openalex.org/W2069087780", "https://openalex.org/W2366695111", "https://openalex.org/W2979561317", "https://openalex.org/W2987737929", "https://openalex.org/W3118321514", "https://openalex.org/W2019985848", "https://openalex.org/W2356158448")'

$ws.Range("F4").Value = 'Hepatology Communications'

$ws.Range("G4").Value = 'https://openalex.org/S4210228550'

$ws.Range("H4").Value = 'Lippincott Williams & Wilkins'

$ws.Range("I4").Value = '2471-254X'

$ws.Range("AF5").Value = 'c("https://openalex.org/W2186218760", "https://openalex.org/W2141629007", "https://openalex.org/W2998456125", "https://openalex.org/W2088646614", "https://openalex.org/W2473498105", "https://openalex.org/W2031943657", "https://openalex.org/W2919691020", "https://openalex.org/W1997688092", "https://openalex.org/W2769989561", "https://openalex.org/W2374585174")'

$ws.Range("AF7").Value = 'c("https://openalex.org/W2417403314", "https://openalex.org/W1640165957", "https://openalex.org/W3105232760", "https://openalex.org/W2461769281", "https://openalex.org/W2412718868", "https://openalex.org/W167219501", "https://openalex.org/W2195790131", "https://openalex.org/W131384799", "https://openalex.org/W4391141357", "https://openalex.org/W2085904991")'

$ws.Range("F11").Value = 'The American Journal of Gastroenterology'

$ws.Range("G11").Value = 'https://openalex.org/S66441642'

$ws.Range("H11").Value = 'Lippincott Williams & Wilkins'

$ws.Range("I11").Value = '0002-9270'

$ws.Range("V11").Value = "'FALSE"
